# Update "想去人数" (column F) counts on the "展览" and "全部类型" sheets,
# matching the regenerated data snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll = $wb.Worksheets.Item("全部类型")

# Sheet "展览": row -> new F value
$exhibitUpdates = @{
    2  = 208
    3  = 5529
    4  = 38
    7  = 657
    8  = 642
    10 = 1078
    12 = 1540
    13 = 5077
    16 = 208
    17 = 19
    18 = 7
    19 = 109
    20 = 4356
    21 = 205
    22 = 1153
    24 = 59
    26 = 56
    27 = 165
    28 = 62
    29 = 147
    31 = 342
    32 = 6
    33 = 40
    35 = 29
    36 = 42
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

# Sheet "全部类型": row -> new F value
$allUpdates = @{
    2  = 209
    4  = 5529
    5  = 38
    8  = 657
    9  = 642
    11 = 1078
    13 = 1540
    14 = 5077
    17 = 208
    18 = 19
    19 = 7
    20 = 109
    21 = 4356
    22 = 205
    23 = 1153
    25 = 59
    27 = 56
    28 = 165
    29 = 62
    30 = 147
    32 = 342
    33 = 6
    34 = 40
    36 = 29
    37 = 42
}

foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}
